$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 76.16539233333333
$ws.Range("H2").Value = 228.496177
$ws.Range("I2").Value = 0.5742142031125765
$ws.Range("J2").Value = 0.5742142031125764
$ws.Range("M2").Value = 39.02618766666667
$ws.Range("N2").Value = 117.078563
$ws.Range("O2").Value = 0.4958819606525626
$ws.Range("P2").Value = 0.4958819606525626
$ws.Range("Q2").Value = 2972.444894905961
$ws.Range("R2").Value = 26752.00405415365
$ws.Range("S2").Value = 0.2847424648740132
$ws.Range("T2").Value = 0.2847424648740132

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 76.16539233333333
$ws.Range("H3").Value = 228.496177
$ws.Range("I3").Value = 0.5742142031125765
$ws.Range("J3").Value = 0.5742142031125764
$ws.Range("O3").Value = 0.03275155884322009
$ws.Range("P3").Value = 0.03275155884322008
$ws.Range("Q3").Value = 196.3213256550607
$ws.Range("R3").Value = 1766.891930895546
$ws.Range("S3").Value = 0.01880641026185428
$ws.Range("T3").Value = 0.01880641026185428

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 76.16539233333333
$ws.Range("H4").Value = 228.496177
$ws.Range("I4").Value = 0.5742142031125765
$ws.Range("J4").Value = 0.5742142031125764
$ws.Range("M4").Value = 35.04673133333333
$ws.Range("N4").Value = 105.140194
$ws.Range("O4").Value = 0.4453174364986936
$ws.Range("P4").Value = 0.4453174364986936
$ws.Range("Q4").Value = 2669.34804200426
$ws.Range("R4").Value = 24024.13237803834
$ws.Range("S4").Value = 0.2557075969312327
$ws.Range("T4").Value = 0.2557075969312327

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 76.16539233333333
$ws.Range("H5").Value = 228.496177
$ws.Range("I5").Value = 0.5742142031125765
$ws.Range("J5").Value = 0.5742142031125764
$ws.Range("M5").Value = 2.050074333333333
$ws.Range("N5").Value = 6.150223
$ws.Range("O5").Value = 0.02604904400552376
$ws.Range("P5").Value = 0.02604904400552376
$ws.Range("Q5").Value = 156.1447159108301
$ws.Range("R5").Value = 1405.302443197471
$ws.Range("S5").Value = 0.01495773104547626
$ws.Range("T5").Value = 0.01495773104547626

$ws.Range("I6").Value = 0.03912478832313545
$ws.Range("J6").Value = 0.03912478832313544
$ws.Range("M6").Value = 39.02618766666667
$ws.Range("N6").Value = 117.078563
$ws.Range("O6").Value = 0.4958819606525626
$ws.Range("P6").Value = 0.4958819606525626
$ws.Range("Q6").Value = 202.5311750997912
$ws.Range("R6").Value = 1822.780575898121
$ws.Range("S6").Value = 0.01940127674379289
$ws.Range("T6").Value = 0.01940127674379289

$ws.Range("I7").Value = 0.03912478832313545
$ws.Range("J7").Value = 0.03912478832313544
$ws.Range("O7").Value = 0.03275155884322009
$ws.Range("P7").Value = 0.03275155884322008
$ws.Range("S7").Value = 0.001281397806993701
$ws.Range("T7").Value = 0.0012813978069937

$ws.Range("I8").Value = 0.03912478832313545
$ws.Range("J8").Value = 0.03912478832313544
$ws.Range("M8").Value = 35.04673133333333
$ws.Range("N8").Value = 105.140194
$ws.Range("O8").Value = 0.4453174364986936
$ws.Range("P8").Value = 0.4453174364986936
$ws.Range("Q8").Value = 181.8792996377998
$ws.Range("R8").Value = 1636.913696740198
$ws.Range("S8").Value = 0.0174229504396127
$ws.Range("T8").Value = 0.0174229504396127

$ws.Range("I9").Value = 0.03912478832313545
$ws.Range("J9").Value = 0.03912478832313544
$ws.Range("M9").Value = 2.050074333333333
$ws.Range("N9").Value = 6.150223
$ws.Range("O9").Value = 0.02604904400552376
$ws.Range("P9").Value = 0.02604904400552376
$ws.Range("Q9").Value = 10.63911154526011
$ws.Range("R9").Value = 95.75200390734101
$ws.Range("S9").Value = 0.001019163332736157
$ws.Range("T9").Value = 0.001019163332736157

$ws.Range("G10").Value = 51.18420533333333
$ws.Range("H10").Value = 153.552616
$ws.Range("I10").Value = 0.3858799485835225
$ws.Range("J10").Value = 0.3858799485835225
$ws.Range("M10").Value = 39.02618766666667
$ws.Range("N10").Value = 117.078563
$ws.Range("O10").Value = 0.4958819606525626
$ws.Range("P10").Value = 0.4958819606525626
$ws.Range("Q10").Value = 1997.524402907867
$ws.Range("R10").Value = 17977.71962617081
$ws.Range("S10").Value = 0.1913509054801072
$ws.Range("T10").Value = 0.1913509054801072

$ws.Range("G11").Value = 51.18420533333333
$ws.Range("H11").Value = 153.552616
$ws.Range("I11").Value = 0.3858799485835225
$ws.Range("J11").Value = 0.3858799485835225
$ws.Range("O11").Value = 0.03275155884322009
$ws.Range("P11").Value = 0.03275155884322008
$ws.Range("Q11").Value = 131.9306674042187
$ws.Range("R11").Value = 1187.376006637968
$ws.Range("S11").Value = 0.01263816984245198
$ws.Range("T11").Value = 0.01263816984245198

$ws.Range("G12").Value = 51.18420533333333
$ws.Range("H12").Value = 153.552616
$ws.Range("I12").Value = 0.3858799485835225
$ws.Range("J12").Value = 0.3858799485835225
$ws.Range("M12").Value = 35.04673133333333
$ws.Range("N12").Value = 105.140194
$ws.Range("O12").Value = 0.4453174364986936
$ws.Range("P12").Value = 0.4453174364986936
$ws.Range("Q12").Value = 1793.8390928275
$ws.Range("R12").Value = 16144.55183544751
$ws.Range("S12").Value = 0.171839069499462
$ws.Range("T12").Value = 0.171839069499462

$ws.Range("G13").Value = 51.18420533333333
$ws.Range("H13").Value = 153.552616
$ws.Range("I13").Value = 0.3858799485835225
$ws.Range("J13").Value = 0.3858799485835225
$ws.Range("M13").Value = 2.050074333333333
$ws.Range("N13").Value = 6.150223
$ws.Range("O13").Value = 0.02604904400552376
$ws.Range("P13").Value = 0.02604904400552376
$ws.Range("Q13").Value = 104.9314256259298
$ws.Range("R13").Value = 944.382830633368
$ws.Range("S13").Value = 0.01005180376150142
$ws.Range("T13").Value = 0.01005180376150142

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.103602
$ws.Range("H14").Value = 0.310806
$ws.Range("I14").Value = 0.0007810599807654878
$ws.Range("J14").Value = 0.0007810599807654877
$ws.Range("M14").Value = 39.02618766666667
$ws.Range("N14").Value = 117.078563
$ws.Range("O14").Value = 0.4958819606525626
$ws.Range("P14").Value = 0.4958819606525626
$ws.Range("Q14").Value = 4.043191094642
$ws.Range("R14").Value = 36.38871985177801
$ws.Range("S14").Value = 0.0003873135546492429
$ws.Range("T14").Value = 0.0003873135546492428

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.103602
$ws.Range("H15").Value = 0.310806
$ws.Range("I15").Value = 0.0007810599807654878
$ws.Range("J15").Value = 0.0007810599807654877
$ws.Range("O15").Value = 0.03275155884322009
$ws.Range("P15").Value = 0.03275155884322008
$ws.Range("Q15").Value = 0.2670409927320001
$ws.Range("R15").Value = 2.403368934588001
$ws.Range("S15").Value = 0.00002558093192012523
$ws.Range("T15").Value = 0.00002558093192012522

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.103602
$ws.Range("H16").Value = 0.310806
$ws.Range("I16").Value = 0.0007810599807654878
$ws.Range("J16").Value = 0.0007810599807654877
$ws.Range("M16").Value = 35.04673133333333
$ws.Range("N16").Value = 105.140194
$ws.Range("O16").Value = 0.4453174364986936
$ws.Range("P16").Value = 0.4453174364986936
$ws.Range("Q16").Value = 3.630911459596001
$ws.Range("R16").Value = 32.678203136364
$ws.Range("S16").Value = 0.0003478196283862059
$ws.Range("T16").Value = 0.0003478196283862059

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.103602
$ws.Range("H17").Value = 0.310806
$ws.Range("I17").Value = 0.0007810599807654878
$ws.Range("J17").Value = 0.0007810599807654877
$ws.Range("M17").Value = 2.050074333333333
$ws.Range("N17").Value = 6.150223
$ws.Range("O17").Value = 0.02604904400552376
$ws.Range("P17").Value = 0.02604904400552376
$ws.Range("Q17").Value = 0.212391801082
$ws.Range("R17").Value = 1.911526209738
$ws.Range("S17").Value = 0.00002034586580991373
$ws.Range("T17").Value = 0.00002034586580991373
